# Natmi following Dr Hou advice
# Update the Pgf-Nrp1 LR-pairs sheet: recompute the edge-expression rows
# for the three clusters (ECs, FAPs, sCs), including three additional
# target-cluster rows that weren't present before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one full row (columns A..T). The leading "," on each
# row forces PowerShell to keep it as a nested array element instead of
# splicing/flattening it into $rowData.
$rowData = @(
  ,@("ECs",  "Pgf", "Nrp1", "ECs",  2, 0.6666666666666666, 2.384646333333333, 7.153938999999999,  0.3245205637288701, 0.3245205637288701, 3, 1, 102.8289443333334, 308.486833,  0.5559120396302444, 0.5559120396302443, 245.2106650650208,  2206.895985585187,  0.1804048884844729,  0.1804048884844729)
  ,@("ECs",  "Pgf", "Nrp1", "FAPs", 2, 0.6666666666666666, 2.384646333333333, 7.153938999999999,  0.3245205637288701, 0.3245205637288701, 3, 1, 63.66262833333334, 190.987885,  0.3441717873742006, 0.3441717873742006, 151.8128532254461,  1366.315679029015,  0.1116908224582484,  0.1116908224582484)
  ,@("ECs",  "Pgf", "Nrp1", "sCs",  2, 0.6666666666666666, 2.384646333333333, 7.153938999999999,  0.3245205637288701, 0.3245205637288701, 3, 1, 18.481835,          55.445505,   0.09991617299555507,0.09991617299555505,44.07264006602166,  396.6537605941949,  0.03242485278614883, 0.03242485278614884)
  ,@("FAPs", "Pgf", "Nrp1", "ECs",  3, 1,                   3.736173,          11.208519,           0.5084464522895362, 0.5084464522895361, 3, 1, 102.8289443333334, 308.486833,  0.5559120396302444, 0.5559120396302443, 384.1867254367031,  3457.680528930328,  0.2826515043350378,  0.2826515043350377)
  ,@("FAPs", "Pgf", "Nrp1", "FAPs", 3, 1,                   3.736173,          11.208519,           0.5084464522895362, 0.5084464522895361, 3, 1, 63.66262833333334, 190.987885,  0.3441717873742006, 0.3441717873742006, 237.854593088035,   2140.691337792315,  0.1749929242685609,  0.1749929242685609)
  ,@("FAPs", "Pgf", "Nrp1", "sCs",  3, 1,                   3.736173,          11.208519,           0.5084464522895362, 0.5084464522895361, 3, 1, 18.481835,          55.445505,   0.09991617299555507,0.09991617299555505,69.05133291745501,  621.4619962570951,  0.05080202368593754, 0.05080202368593752)
  ,@("sCs",  "Pgf", "Nrp1", "ECs",  3, 1,                   1.227394,          3.682182,            0.1670329839815937, 0.1670329839815937, 3, 1, 102.8289443333334, 308.486833,  0.5559120396302444, 0.5559120396302443, 126.2116293010674,  1135.904663709606,  0.0928556468107337,  0.09285564681073369)
  ,@("sCs",  "Pgf", "Nrp1", "FAPs", 3, 1,                   1.227394,          3.682182,            0.1670329839815937, 0.1670329839815937, 3, 1, 63.66262833333334, 190.987885,  0.3441717873742006, 0.3441717873742006, 78.13912804056335,  703.2521523650701,  0.05748804064739133, 0.05748804064739133)
  ,@("sCs",  "Pgf", "Nrp1", "sCs",  3, 1,                   1.227394,          3.682182,            0.1670329839815937, 0.1670329839815937, 3, 1, 18.481835,          55.445505,   0.09991617299555507,0.09991617299555505,22.68449338799,     204.16044049191,    0.01668929652346869, 0.01668929652346869)
)

$startRow = 2
for ($r = 0; $r -lt $rowData.Length; $r++) {
  $vals = $rowData[$r]
  $rn = $startRow + $r
  for ($c = 0; $c -lt $vals.Length; $c++) {
    $ws.Cells.Item($rn, $c + 1).Value = $vals[$c]
  }
}
